# #1189143 - changed the conference numbers to reflect the new format - with Michael Adams
#
# The "Conference" column (B) on the main schedule sheet used plain numeric
# codes (11, 12, 13, 21, 22, 23). They are replaced with the new
# alphanumeric conference naming scheme: 1a, 1b, 1c, 2a, 2b, 2c.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B3:B9").Value   = "1a"
$ws.Range("B10:B14").Value = "1b"
$ws.Range("B15:B20").Value = "1c"
$ws.Range("B21:B28").Value = "2a"
$ws.Range("B29:B36").Value = "2b"
$ws.Range("B37:B42").Value = "2c"
